# Trade #26 closed/recorded at 2026-02-16 22:54:40 - base_strategy DOWN +0.000%
# Appends the new trade row (row 27) to both the "All Trades" log sheet and
# the per-strategy "base_strategy" sheet, mirroring the existing row layout.

$wb = $excel.ActiveWorkbook

$rowNum      = 27
$tradeNum    = 26
$tradeDate   = "2026-02-16"
$tradeTime   = "22:54:40"
$strategy    = "base_strategy"
$side        = "DOWN"
$entryPrice  = 49.999998
$status      = "OPEN"
$pnlPct      = 0
$pnlUsd      = 0
$capitalAfter = 100
$entrySlippage = 0
$exitSlippage  = 0
$confidence  = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$duration    = 0

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($rowNum, 1).Value = $tradeNum

    # Column B holds a plain "yyyy-mm-dd" text label (not a real date), so
    # force the cell to Text first - otherwise Excel auto-parses the
    # ISO-looking string into a date serial number.
    $dateCell = $ws.Cells.Item($rowNum, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $tradeDate

    $ws.Cells.Item($rowNum, 3).Value = $tradeTime
    $ws.Cells.Item($rowNum, 4).Value = $strategy
    $ws.Cells.Item($rowNum, 5).Value = $side
    $ws.Cells.Item($rowNum, 6).Value = $entryPrice

    # Column G (Exit Price) stays blank - trade is still OPEN. A lone
    # apostrophe is Excel's own "empty text" quote-prefix marker, so the
    # cell keeps text type but reads as an empty string (matching the
    # source row layout) instead of collapsing to a truly empty/number cell.
    $ws.Cells.Item($rowNum, 7).Value = "'"

    $ws.Cells.Item($rowNum, 8).Value = $status
    $ws.Cells.Item($rowNum, 9).Value = $pnlPct
    $ws.Cells.Item($rowNum, 10).Value = $pnlUsd
    $ws.Cells.Item($rowNum, 11).Value = $capitalAfter
    $ws.Cells.Item($rowNum, 12).Value = $entrySlippage
    $ws.Cells.Item($rowNum, 13).Value = $exitSlippage
    $ws.Cells.Item($rowNum, 14).Value = $confidence
    $ws.Cells.Item($rowNum, 15).Value = $entryReason

    # Column P (Exit Reason) stays blank - trade is still OPEN. Same
    # empty-text trick as column G above.
    $ws.Cells.Item($rowNum, 16).Value = "'"

    $ws.Cells.Item($rowNum, 17).Value = $duration
}
